# Update "how often" (clean_window) specifications on the "outcomes" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

# Update clean_window (column C) values for the existing rows whose
# specification changed but whose cohort id / name stayed the same.
$ws.Cells.Item(2, 3).Value = 30
$ws.Cells.Item(7, 3).Value = 365
$ws.Cells.Item(18, 3).Value = 90
$ws.Cells.Item(26, 3).Value = 9999
$ws.Cells.Item(32, 3).Value = 9999
$ws.Cells.Item(39, 3).Value = 180
$ws.Cells.Item(47, 3).Value = 180
$ws.Cells.Item(49, 3).Value = 365
$ws.Cells.Item(52, 3).Value = 365
$ws.Cells.Item(54, 3).Value = 30
$ws.Cells.Item(92, 3).Value = 180

# Remove the row for cohort_definition_id 1017 ("Earliest event of Neonatal
# Thrombocytopenia (NT), less than 1 year old"). This shifts every
# subsequent row up by one, shrinking the used range from C263 to C262.
$ws.Rows.Item(214).Delete()
